# Adds a new "2022-Q3" worksheet (with the latest quarterly fund-holding
# data) right after the "总计" summary sheet and before "2022-Q2", and
# updates the "总计" summary sheet with a new row for 2022-Q3 (the other
# quarterly sheets keep their own data untouched - they simply shift one
# tab position to the right because of the insertion).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) summary sheet: add 2022-Q3 as the newest
#    row right after the header, pushing the existing rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 9, 0.42),
    @(1, "2022-Q2", 7, 0.47),
    @(2, "2022-Q1", 7, 0.61),
    @(3, "2021-Q3", 7, 0.61),
    @(4, "2021-Q2", 8, 0.63),
    @(5, "2020-Q4", 8, 0.32)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# Match the existing style used on column A (bold / centred / bordered)
# for the newly-added row 7.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# ---------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q3" worksheet right before "2022-Q2".
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($existingQ2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3Rows = @(
    @(0, "005613", "上投摩根富时发达市场REITs指数（QDII）人民币份额", "3.41", "92.80", "2.98", "0.1016", 9),
    @(1, "005614", "上投摩根富时发达市场REITs指数（QDII）美钞",       "3.41", "92.80", "2.98", "0.1016", 9),
    @(2, "005615", "上投摩根富时发达市场REITs指数（QDII）美汇",       "3.41", "92.80", "2.98", "0.1016", 9),
    @(3, "000179", "广发美国房地产指数（QDII）人民币A",               "1.82", "92.37", "2.42", "0.0440", 10),
    @(4, "000180", "广发美国房地产指数（QDII）美元A",                 "1.82", "92.37", "2.42", "0.0440", 10),
    @(5, "160140", "南方道琼斯美国精选REIT指数（QDII-LOF）A",          "0.78", "91.13", "2.45", "0.0191", 10),
    @(6, "160141", "南方道琼斯美国精选REIT指数（QDII-LOF）C",          "0.42", "91.13", "2.45", "0.0103", 10),
    @(7, "016278", "广发美国房地产指数（QDII）人民币C",               "0.01", "92.37", "2.42", "0.0002", 10),
    @(8, "016279", "广发美国房地产指数（QDII）美元C",                 "0.01", "92.37", "2.42", "0.0002", 10)
)

for ($i = 0; $i -lt $q3Rows.Length; $i++) {
    $r = $i + 2
    $row = $q3Rows[$i]

    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]

    # D/E/F/G are stored as text in the source data, not as numbers -
    # force text formatting before assigning so "92.80"/"0.0002" etc.
    # keep their exact original representation.
    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[6]

    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Copy the header / column-A formatting from the neighbouring "2022-Q2"
# sheet so the new tab matches the existing look (bold, centred, thin
# border on header row and row-index column).
$existingQ2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$existingQ2.Range("A2").Copy()
$q3.Range("A2:A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$summary.Activate()
